$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 99.75
$ws.Range("I12").Value = 99.666664
$ws.Range("K12").Value = 99.666664
$ws.Range("M12").Value = 70.333336
$ws.Range("H32").Value = 962.75
$ws.Range("I32").Value = 1200
$ws.Range("J32").Value = 941.1818
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 941.1818
$ws.Range("M32").Value = -874
$ws.Range("N32").Value = -1593.1818
$ws.Range("H76").Value = 2472629.2
$ws.Range("I76").Value = 3370767.2
$ws.Range("J76").Value = 2750
$ws.Range("K76").Value = 3370767.2
$ws.Range("L76").Value = 2750
$ws.Range("M76").Value = -3370452.2
$ws.Range("N76").Value = -3380
$ws.Range("H79").Value = 2472629.2
$ws.Range("I79").Value = 3370767.2
$ws.Range("J79").Value = 2750
$ws.Range("K79").Value = 3370767.2
$ws.Range("L79").Value = 2750
$ws.Range("M79").Value = -3369675.2
$ws.Range("N79").Value = -4934
$ws.Range("H98").Value = 4368458
$ws.Range("I98").Value = 5879265.5
$ws.Range("J98").Value = 3904.111
$ws.Range("K98").Value = 5879265.5
$ws.Range("L98").Value = 3904.111
$ws.Range("M98").Value = -5877767.5
$ws.Range("N98").Value = -6900.111
$ws.Range("H115").Value = 948.9286
$ws.Range("I115").Value = 496.25
$ws.Range("J115").Value = 1130
$ws.Range("K115").Value = 1488.75
$ws.Range("L115").Value = 3390
$ws.Range("M115").Value = 78.25
$ws.Range("N115").Value = -6524
$ws.Range("H116").Value = 4118.857
$ws.Range("I116").Value = 9125
$ws.Range("K116").Value = 9125
$ws.Range("M116").Value = -5683
$ws.Range("H122").Value = 4368458
$ws.Range("I122").Value = 5879265.5
$ws.Range("J122").Value = 3904.111
$ws.Range("K122").Value = 17637796.5
$ws.Range("L122").Value = 11712.333
$ws.Range("M122").Value = -17635346.5
$ws.Range("N122").Value = -16612.333
$ws.Range("H125").Value = 1143
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1143
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 10287
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -15207
$ws.Range("H138").Value = 2395.6382
$ws.Range("I138").Value = 908.25
$ws.Range("J138").Value = 5568.7334
$ws.Range("K138").Value = 2724.75
$ws.Range("L138").Value = 16706.2002
$ws.Range("M138").Value = 2415.25
$ws.Range("N138").Value = -26986.2002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 765.1739
$ws.Range("I74").Value = 732.3333
$ws.Range("J74").Value = 776.7646999999999
$ws.Range("K74").Value = 732.3333
$ws.Range("L74").Value = 776.7646999999999
$ws.Range("M74").Value = 141.6667
$ws.Range("N74").Value = -2524.7647
$ws.Range("H77").Value = 765.1739
$ws.Range("I77").Value = 732.3333
$ws.Range("J77").Value = 776.7646999999999
$ws.Range("K77").Value = 3661.6665
$ws.Range("L77").Value = 3883.8235
$ws.Range("M77").Value = 706.3334999999997
$ws.Range("N77").Value = -12619.8235

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 74.07143000000001
$ws.Range("I7").Value = 52
$ws.Range("J7").Value = 90.625
$ws.Range("K7").Value = 52
$ws.Range("L7").Value = 90.625
$ws.Range("M7").Value = 61
$ws.Range("N7").Value = -316.625
$ws.Range("H129").Value = 31793
$ws.Range("I129").Value = 9999
$ws.Range("J129").Value = 49954.668
$ws.Range("K129").Value = 9999
$ws.Range("L129").Value = 49954.668
$ws.Range("M129").Value = -4999
$ws.Range("N129").Value = -59954.668
$ws.Range("H137").Value = 30000
$ws.Range("I137").Value = 30000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 30000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -24900
$ws.Range("N137").Value = ""

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 304.63635
$ws.Range("J2").Value = 343.21054
$ws.Range("L2").Value = 2059.26324
$ws.Range("N2").Value = -2285.26324
$ws.Range("H4").Value = 1278.125
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 1818.1818
$ws.Range("K4").Value = 270
$ws.Range("L4").Value = 5454.5454
$ws.Range("M4").Value = -158
$ws.Range("N4").Value = -5678.5454
$ws.Range("H6").Value = 324.375
$ws.Range("I6").Value = 148.33333
$ws.Range("J6").Value = 430
$ws.Range("K6").Value = 444.99999
$ws.Range("L6").Value = 1290
$ws.Range("M6").Value = -331.99999
$ws.Range("N6").Value = -1516
$ws.Range("H12").Value = 253.09091
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 242.66667
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 728.00001
$ws.Range("M12").Value = -727
$ws.Range("N12").Value = -1074.00001
$ws.Range("H104").Value = 3700
$ws.Range("I104").Value = 3500
$ws.Range("J104").Value = 3900
$ws.Range("K104").Value = 10500
$ws.Range("L104").Value = 11700
$ws.Range("M104").Value = -7879
$ws.Range("N104").Value = -16942
$ws.Range("H111").Value = 50000492
$ws.Range("I111").Value = 50000492
$ws.Range("K111").Value = 150001476
$ws.Range("M111").Value = -149998409

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1264.8485
$ws.Range("I97").Value = 957.7778
$ws.Range("J97").Value = 2646.6667
$ws.Range("K97").Value = 957.7778
$ws.Range("L97").Value = 2646.6667
$ws.Range("M97").Value = -461.7778
$ws.Range("N97").Value = -3638.6667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 65753.375
$ws.Range("I136").Value = 100505.4
$ws.Range("J136").Value = 7833.3335
$ws.Range("K136").Value = 301516.2
$ws.Range("L136").Value = 23500.0005
$ws.Range("M136").Value = -298966.2
$ws.Range("N136").Value = -28600.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 32259930
$ws.Range("I132").Value = 43479948
$ws.Range("J132").Value = 2372.625
$ws.Range("K132").Value = 130439844
$ws.Range("L132").Value = 7117.875
$ws.Range("M132").Value = -130437314
$ws.Range("N132").Value = -12177.875
$ws.Range("H136").Value = 23742.395
$ws.Range("I136").Value = 27306.162
$ws.Range("J136").Value = 1765.8334
$ws.Range("K136").Value = 81918.486
$ws.Range("L136").Value = 5297.5002
$ws.Range("M136").Value = -79368.486
$ws.Range("N136").Value = -10397.5002
